$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.6609919999999999
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 2.110264333333333
$ws.Range("N2").Value = 6.330793
$ws.Range("O2").Value = 0.3832041185227171
$ws.Range("P2").Value = 0.3832041185227171
$ws.Range("Q2").Value = 1.394867842218666
$ws.Range("R2").Value = 12.553810579968
$ws.Range("S2").Value = 0.3832041185227171
$ws.Range("T2").Value = 0.3832041185227171

# Row 3
$ws.Range("G3").Value = 0.6609919999999999
$ws.Range("O3").Value = 0.0946183755984393
$ws.Range("P3").Value = 0.0946183755984393
$ws.Range("R3").Value = 3.099708764159999
$ws.Range("S3").Value = 0.0946183755984393
$ws.Range("T3").Value = 0.0946183755984393

# Row 4
$ws.Range("G4").Value = 0.6609919999999999
$ws.Range("M4").Value = 2.065388333333333
$ws.Range("N4").Value = 6.196165
$ws.Range("O4").Value = 0.3750550597762889
$ws.Range("P4").Value = 0.3750550597762889
$ws.Range("Q4").Value = 1.365205165226667
$ws.Range("R4").Value = 12.28684648704
$ws.Range("S4").Value = 0.3750550597762889
$ws.Range("T4").Value = 0.3750550597762889

# Row 5
$ws.Range("G5").Value = 0.6609919999999999
$ws.Range("M5").Value = 0.8101876666666666
$ws.Range("N5").Value = 2.430563
$ws.Range("O5").Value = 0.1471224461025547
$ws.Range("P5").Value = 0.1471224461025547
$ws.Range("Q5").Value = 0.5355275661653333
$ws.Range("R5").Value = 4.819748095487999
$ws.Range("S5").Value = 0.1471224461025547
$ws.Range("T5").Value = 0.1471224461025547
